$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Name" column header in G3, matching the bold/centered header style
# already used by A3 ("Part").
$ws.Range("G3").Value = "Name"
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").HorizontalAlignment = -4108

# New "Name" values for rows 4-12 (entered in this order so the shared
# string table lines up with how the workbook was authored), formatted as
# text like the "Manufacturer part number" column already on each row.
$ws.Range("G6").Value = "C001-256"
$ws.Range("G9").Value = "LED001-256"
$ws.Range("G11").Value = "R001-256"
$ws.Range("G7").Value = "C257-288"
$ws.Range("G12").Value = "IC001-256"
$ws.Range("G8").Value = "IC257-288"
$ws.Range("G10").Value = "R257-557"
$ws.Range("G5").Value = "J3-4"
$ws.Range("G4").Value = "J1-2"

$ws.Range("G4:G12").NumberFormat = "@"

# Selection, as left by the editing session.
$ws.Range("H5").Select()
